# Apply scheduled market-data update to Chocobo Profits workbook
$wb = $excel.ActiveWorkbook

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3059.158
$ws.Range("I98").Value = 1455.2858
$ws.Range("J98").Value = 7550
$ws.Range("K98").Value = 1455.2858
$ws.Range("L98").Value = 7550
$ws.Range("M98").Value = 42.71419999999989
$ws.Range("N98").Value = -10546

# ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 938.4167
$ws.Range("I111").Value = 903.8570999999999
$ws.Range("J111").Value = 986.8
$ws.Range("K111").Value = 2711.5713
$ws.Range("L111").Value = 2960.4
$ws.Range("M111").Value = 355.4287000000004
$ws.Range("N111").Value = -9094.4

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4932
$ws.Range("I113").Value = 1705
$ws.Range("J113").Value = 5162.5
$ws.Range("K113").Value = 1705
$ws.Range("L113").Value = 5162.5
$ws.Range("M113").Value = 1549
$ws.Range("N113").Value = -11670.5

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 329311.97
$ws.Range("I116").Value = 836322.0600000001
$ws.Range("J116").Value = 9095.053
$ws.Range("K116").Value = 836322.0600000001
$ws.Range("L116").Value = 9095.053
$ws.Range("M116").Value = -832880.0600000001
$ws.Range("N116").Value = -15979.053

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 3059.158
$ws.Range("I122").Value = 1455.2858
$ws.Range("J122").Value = 7550
$ws.Range("K122").Value = 4365.857400000001
$ws.Range("L122").Value = 22650
$ws.Range("M122").Value = -1915.857400000001
$ws.Range("N122").Value = -27550

# ALC row 123
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = $null

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4376.58
$ws.Range("J138").Value = 5450.2856
$ws.Range("L138").Value = 16350.8568
$ws.Range("N138").Value = -26630.8568

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1177.196
$ws.Range("I61").Value = 1161.1
$ws.Range("J61").Value = 1235.7273
$ws.Range("K61").Value = 1161.1
$ws.Range("L61").Value = 1235.7273
$ws.Range("M61").Value = -949.0999999999999
$ws.Range("N61").Value = -1659.7273

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3138.2327
$ws.Range("I74").Value = 3896.5833
$ws.Range("J74").Value = 2180.3157
$ws.Range("K74").Value = 3896.5833
$ws.Range("L74").Value = 2180.3157
$ws.Range("M74").Value = -3022.5833
$ws.Range("N74").Value = -3928.3157

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3138.2327
$ws.Range("I77").Value = 3896.5833
$ws.Range("J77").Value = 2180.3157
$ws.Range("K77").Value = 19482.9165
$ws.Range("L77").Value = 10901.5785
$ws.Range("M77").Value = -15114.9165
$ws.Range("N77").Value = -19637.5785

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1177.196
$ws.Range("I136").Value = 1161.1
$ws.Range("J136").Value = 1235.7273
$ws.Range("K136").Value = 3483.3
$ws.Range("L136").Value = 3707.1819
$ws.Range("M136").Value = -933.2999999999997
$ws.Range("N136").Value = -8807.1819

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1242
$ws.Range("I107").Value = 1242
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1242
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 678
$ws.Range("N107").Value = $null

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1921.9048
$ws.Range("I31").Value = 755.89655
$ws.Range("J31").Value = 4523
$ws.Range("K31").Value = 755.89655
$ws.Range("L31").Value = 4523
$ws.Range("M31").Value = -460.89655
$ws.Range("N31").Value = -5113

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1921.9048
$ws.Range("I34").Value = 755.89655
$ws.Range("J34").Value = 4523
$ws.Range("K34").Value = 755.89655
$ws.Range("L34").Value = 4523
$ws.Range("M34").Value = -553.89655
$ws.Range("N34").Value = -4927

# CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 7666.3335
$ws.Range("I86").Value = 6500
$ws.Range("K86").Value = 6500
$ws.Range("M86").Value = -5377

# CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 7666.3335
$ws.Range("I89").Value = 6500
$ws.Range("K89").Value = 32500
$ws.Range("M89").Value = -26884

# CUL row 3
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2362.3076

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 760.25
$ws.Range("I131").Value = 549
$ws.Range("J131").Value = 764.5612
$ws.Range("K131").Value = 1647
$ws.Range("L131").Value = 2293.6836
$ws.Range("M131").Value = 3393
$ws.Range("N131").Value = -12373.6836

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4351.2
$ws.Range("I122").Value = 3126.5
$ws.Range("J122").Value = 9250
$ws.Range("K122").Value = 9379.5
$ws.Range("L122").Value = 27750
$ws.Range("M122").Value = -6929.5
$ws.Range("N122").Value = -32650

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2414.6
$ws.Range("I132").Value = 1837.7894
$ws.Range("J132").Value = 3099.5625
$ws.Range("K132").Value = 5513.3682
$ws.Range("L132").Value = 9298.6875
$ws.Range("M132").Value = -2983.3682
$ws.Range("N132").Value = -14358.6875

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 483.70834
$ws.Range("I16").Value = 483.70834
$ws.Range("K16").Value = 483.70834
$ws.Range("M16").Value = -313.70834

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7799.3335
$ws.Range("I40").Value = 6743.5557
$ws.Range("J40").Value = 10966.667
$ws.Range("K40").Value = 6743.5557
$ws.Range("L40").Value = 10966.667
$ws.Range("M40").Value = -6607.5557
$ws.Range("N40").Value = -11238.667

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 716.42426
$ws.Range("I68").Value = 716.42426
$ws.Range("K68").Value = 716.42426
$ws.Range("M68").Value = 32.57574

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 716.42426
$ws.Range("I71").Value = 716.42426
$ws.Range("K71").Value = 3582.1213
$ws.Range("M71").Value = 161.8787000000002

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4466.6665
$ws.Range("I122").Value = 1950
$ws.Range("K122").Value = 5850
$ws.Range("M122").Value = -3400

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2957.889
$ws.Range("I132").Value = 2158.5454
$ws.Range("J132").Value = 4214
$ws.Range("K132").Value = 6475.6362
$ws.Range("L132").Value = 12642
$ws.Range("M132").Value = -3945.6362
$ws.Range("N132").Value = -17702

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2867.9048
$ws.Range("I136").Value = 1122.6
$ws.Range("J136").Value = 4454.5454
$ws.Range("K136").Value = 3367.8
$ws.Range("L136").Value = 13363.6362
$ws.Range("M136").Value = -817.7999999999997
$ws.Range("N136").Value = -18463.6362

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5012
$ws.Range("I136").Value = 3020.1667
$ws.Range("J136").Value = 7999.75
$ws.Range("K136").Value = 9060.500100000001
$ws.Range("L136").Value = 23999.25
$ws.Range("M136").Value = -6510.500100000001
$ws.Range("N136").Value = -29099.25
